$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(7).Insert()

$ws.Range("A7").Value = 10
$ws.Range("B7").Value = "Vega Modelo de Temuco"
$ws.Range("C7").Value = "La Araucanía"
$ws.Range("D7").Value = "2022-03-10"
$ws.Range("E7").Value = 9
$ws.Range("F7").Value = "Fruta"
$ws.Range("G7").Value = 100108
$ws.Range("H7").Value = "Tropicales y subtropicales"
$ws.Range("I7").Value = 100108004
$ws.Range("J7").Value = "Papaya"
$ws.Range("K7").Value = "Cultivar IV Región"
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 40
$ws.Range("N7").Value = 20000
$ws.Range("O7").Value = 20000
$ws.Range("P7").Value = 20000
$ws.Range("Q7").Value = '$/bandeja 10 kilos'
$ws.Range("R7").Value = "Provincia del Elquí"
$ws.Range("S7").Value = 2000
$ws.Range("T7").Value = 10
